$wb = $excel.ActiveWorkbook

# --- Sheet: classFields ---
# The field/enum-constant rows were reordered (and a couple of
# modifier/type values were swapped along with them) per the fixed
# reflection-derived ordering used by the MSM extractor.
$fields = $wb.Worksheets.Item("classFields")

# domain.Order (rows 2-8)
$fields.Range("B2").Value = "status"
$fields.Range("D2").Value = "domain.OrderStatus"

$fields.Range("B3").Value = "price"
$fields.Range("D3").Value = "int"

$fields.Range("B4").Value = "productId"
$fields.Range("D4").Value = "java.lang.Long"

$fields.Range("B5").Value = "source"
$fields.Range("D5").Value = "domain.OrderSource"

$fields.Range("B6").Value = "productCount"
$fields.Range("D6").Value = "int"

$fields.Range("B7").Value = "customerId"
$fields.Range("D7").Value = "java.lang.Long"

$fields.Range("B8").Value = "id"
$fields.Range("D8").Value = "java.lang.Long"

# domain.OrderStatus (rows 9-15)
$fields.Range("B9").Value = "ROLLBACK"
$fields.Range("C9").Value = "public"
$fields.Range("D9").Value = "domain.OrderStatus"

$fields.Range("B10").Value = "REJECT"
$fields.Range("C10").Value = "public"
$fields.Range("D10").Value = "domain.OrderStatus"

$fields.Range("B11").Value = "ACCEPT"
$fields.Range("C11").Value = "public"
$fields.Range("D11").Value = "domain.OrderStatus"

$fields.Range("B12").Value = "NEW"
$fields.Range("C12").Value = "public"
$fields.Range("D12").Value = "domain.OrderStatus"

$fields.Range("B13").Value = "CONFIRMED"
$fields.Range("C13").Value = "public"
$fields.Range("D13").Value = "domain.OrderStatus"

$fields.Range("B14").Value = "`$VALUES"
$fields.Range("C14").Value = "private"
$fields.Range("D14").Value = "domain.OrderStatus[]"

$fields.Range("B15").Value = "REJECTED"
$fields.Range("C15").Value = "public"
$fields.Range("D15").Value = "domain.OrderStatus"

# domain.Order$OrderBuilder (rows 18-24)
$fields.Range("B18").Value = "status"
$fields.Range("D18").Value = "domain.OrderStatus"

$fields.Range("B19").Value = "price"
$fields.Range("D19").Value = "int"

$fields.Range("B20").Value = "productCount"
$fields.Range("D20").Value = "int"

$fields.Range("B21").Value = "source"
$fields.Range("D21").Value = "domain.OrderSource"

$fields.Range("B22").Value = "id"
$fields.Range("D22").Value = "java.lang.Long"

$fields.Range("B23").Value = "customerId"
$fields.Range("D23").Value = "java.lang.Long"

$fields.Range("B24").Value = "productId"
$fields.Range("D24").Value = "java.lang.Long"

# domain.OrderSource (rows 25-27)
$fields.Range("B25").Value = "PAYMENT"
$fields.Range("C25").Value = "public"
$fields.Range("D25").Value = "domain.OrderSource"

$fields.Range("B26").Value = "`$VALUES"
$fields.Range("C26").Value = "private"
$fields.Range("D26").Value = "domain.OrderSource[]"

$fields.Range("B27").Value = "STOCK"
$fields.Range("C27").Value = "public"
$fields.Range("D27").Value = "domain.OrderSource"

# --- Sheet: methodNumberOfLines ---
$mloc = $wb.Worksheets.Item("methodNumberOfLines")

# domain.OrderStatus synthetic enum methods now report 0 lines
$mloc.Range("C20").Value = 0
$mloc.Range("C21").Value = 0
$mloc.Range("C22").Value = 0
$mloc.Range("C23").Value = 0

# domain.OrderSource synthetic enum methods now report 0 lines
$mloc.Range("C34").Value = 0
$mloc.Range("C35").Value = 0
$mloc.Range("C36").Value = 0
$mloc.Range("C37").Value = 0
